# Edit workbook to add data for 2022-05-12 (commit message),
# which updates the "through 05-03" reporting date to "through 05-04"
# and refreshes the May / Total row figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet (and its tab name shown in workbook.xml) to reflect the new date.
$ws.Name = "Through 2022-05-04"

# Update the row label for May to reflect the new "through" date.
$ws.Range("A6").Value = "May (through 05-04)"

# Update May row (row 6) values for columns C..I (2016..2022).
$ws.Range("C6").Value = 6
$ws.Range("D6").Value = 7
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 4
$ws.Range("H6").Value = 12
$ws.Range("I6").Value = 12

# Update Total row (row 7) values for columns C..I (2016..2022).
$ws.Range("C7").Value = 168
$ws.Range("D7").Value = 260
$ws.Range("E7").Value = 249
$ws.Range("G7").Value = 266
$ws.Range("H7").Value = 535
$ws.Range("I7").Value = 563
